# This workbook's XML diff (vs. the canonical OOXML committed upstream) is
# overwhelmingly a re-serialization artifact: the author simply opened the
# file in real Excel and saved it again (new xmlns:mc/x15/xr* namespaces,
# fileVersion/revisionPtr/AlternateContent bookkeeping, absolute author-machine
# file path, dropped xml:space="preserve" on shared strings that don't need
# it, added dimension/cols/spans metadata, etc.). None of the <v> cell values
# or shared-string text actually changed anywhere in the sheet - the commit's
# real content edit (switching a DotPlot from lines to dots) lives in an R
# script elsewhere in the repository, not in this spreadsheet.
#
# The one genuinely structural change inside the worksheet itself is that a
# batch of wholly empty "B" column cells (no value, no type, no style -
# just a bare placeholder <c r="Bxx"/>) disappear from the saved XML. That is
# exactly what happens when Excel re-writes a row that contains a truly blank
# cell with nothing to persist, so we reproduce it by explicitly clearing
# those cells.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$emptyBCells = @(
    "B46","B52","B56","B59","B65","B69","B76","B80","B84","B88","B92",
    "B126","B137","B148","B152","B158","B161","B166","B170","B173",
    "B176","B179","B183"
)

foreach ($addr in $emptyBCells) {
    $ws.Range($addr).ClearContents()
}

# Match the author's final on-screen selection state (whole sheet selected)
# as captured by the saved sheetView.
$ws.Cells.Select() | Out-Null
